# Weekly update: insert a new week's worth of Acelga price rows
# (Primera / Segunda) at the top of the data block (rows 550-551 on
# Sheet1), pushing all subsequent rows down by two. This mirrors how a
# new week's price observations get prepended to the existing
# chronological log.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 550-551; everything from the old row 550
# onward shifts down to 552 onward (old 584 -> 586, old 585 -> 587).
$ws.Range("A550:A551").EntireRow.Insert()

# New row 550: Acelga / Primera
$ws.Cells.Item(550, 1).Value = 9
$ws.Cells.Item(550, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(550, 3).Value = "Metropolitana"
$ws.Cells.Item(550, 4).Value = 44714
$ws.Cells.Item(550, 5).Value = 13
$ws.Cells.Item(550, 6).Value = 100112009
$ws.Cells.Item(550, 7).Value = "Acelga"
$ws.Cells.Item(550, 8).Value = "Sin especificar"
$ws.Cells.Item(550, 9).Value = "Primera"
$ws.Cells.Item(550, 10).Value = 61
$ws.Cells.Item(550, 11).Value = 11000
$ws.Cells.Item(550, 12).Value = 11000
$ws.Cells.Item(550, 13).Value = 11000
$ws.Cells.Item(550, 14).Value = "`$/docena de atados"
$ws.Cells.Item(550, 15).Value = "Región Metropolitana"
$ws.Cells.Item(550, 16).Value = 3667
$ws.Cells.Item(550, 17).Value = 3
$ws.Cells.Item(550, 18).Value = "Hortaliza"

# New row 551: Acelga / Segunda
$ws.Cells.Item(551, 1).Value = 9
$ws.Cells.Item(551, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(551, 3).Value = "Metropolitana"
$ws.Cells.Item(551, 4).Value = 44714
$ws.Cells.Item(551, 5).Value = 13
$ws.Cells.Item(551, 6).Value = 100112009
$ws.Cells.Item(551, 7).Value = "Acelga"
$ws.Cells.Item(551, 8).Value = "Sin especificar"
$ws.Cells.Item(551, 9).Value = "Segunda"
$ws.Cells.Item(551, 10).Value = 34
$ws.Cells.Item(551, 11).Value = 9000
$ws.Cells.Item(551, 12).Value = 9000
$ws.Cells.Item(551, 13).Value = 9000
$ws.Cells.Item(551, 14).Value = "`$/docena de atados"
$ws.Cells.Item(551, 15).Value = "Región Metropolitana"
$ws.Cells.Item(551, 16).Value = 3000
$ws.Cells.Item(551, 17).Value = 3
$ws.Cells.Item(551, 18).Value = "Hortaliza"
